# language.xlsx tweaks:
#  - insert a new "summaryTitle"/"SUMMARY" row (post-level summary for knowledge)
#  - pretty-up the post game pop-quiz level-select tutorial copy
#  - leave the cursor on row 19 (some more level select tweaks)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the old row 23 ("climate"/"Climate") and populate it.
$ws.Rows(23).Insert() | Out-Null
$ws.Range("A23").Value = "summaryTitle"
$ws.Range("B23").Value = "SUMMARY"

# Tidy up the level-select tutorial copy (now rows 101/102 after the insert).
$ws.Range("B101").Value = "This is a satellite map of Earth. Here you will help us find a sutable location to land our seedling."
$ws.Range("B102").Value = "You can press the image of the requested climate to get more information."

# Match the author's final on-screen selection.
$ws.Activate()
$ws.Rows(19).Select() | Out-Null
